# adding averages and more checks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# --- Row 3 (LOTO (SOPs)): training re-checked -> now expiring sooner / NOT VALID ---
$ws.Range("H3").Value = 20
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "16-Sep-2025"
$ws.Range("J3").Value = "NOT VALID"

# --- Row 4 (Endangered by Electricity A safety Training (SOPs)): refreshed values ---
$ws.Range("H4").Value = -42
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "16-Sep-2025"
$ws.Range("J4").Value = "NOT VALID"

# Row 3 is NOT VALID now too, so give it the same "not valid" pink highlight
# already used on row 4.
$ws.Range("A3:K3").Interior.Color = 13551615
$ws.Range("A4:K4").Interior.Color = 13551615

# Re-balance the title/header font: bold white text (drop the old oversized
# 14pt title face) so the header band reads clearly against its navy fill.
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 11
$ws.Range("A1").Font.Color = 16777215

$ws.Range("A2:K2").Font.Bold = $true
$ws.Range("A2:K2").Font.Size = 11
$ws.Range("A2:K2").Font.Color = 16777215
